$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.082916074980141
$ws.Range("D2").Value = 1.071212724985838
$ws.Range("E2").Value = 1.083512288508778
$ws.Range("F2").Value = 1.08916585565232
$ws.Range("I2").Value = 1.049259449867071
$ws.Range("J2").Value = 1.087783353312988
$ws.Range("K2").Value = 1.073910498376218
$ws.Range("L2").Value = 1.08617771193821
$ws.Range("M2").Value = 1.091816673715599
$ws.Range("N2").Value = 1.089328130445187
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.085554731692022
$ws.Range("D3").Value = 1.07254434795961
$ws.Range("E3").Value = 1.0857416963594
$ws.Range("F3").Value = 1.091161266477799
$ws.Range("I3").Value = 1.049713981405536
$ws.Range("J3").Value = 1.090076559678816
$ws.Range("K3").Value = 1.075057359250172
$ws.Range("L3").Value = 1.088222555074974
$ws.Range("M3").Value = 1.093629163508242
$ws.Range("N3").Value = 1.091624593427089
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.087254164104436
$ws.Range("D4").Value = 1.073400517263951
$ws.Range("E4").Value = 1.087176939231998
$ws.Range("F4").Value = 1.092445515686433
$ws.Range("I4").Value = 1.050003939803736
$ws.Range("J4").Value = 1.091552481088305
$ws.Range("K4").Value = 1.075793439405589
$ws.Range("L4").Value = 1.089538018507015
$ws.Range("M4").Value = 1.094794644488092
$ws.Range("N4").Value = 1.09310261081427
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.087966752775386
$ws.Range("D5").Value = 1.073759160123223
$ws.Range("E5").Value = 1.087778602262117
$ws.Range("F5").Value = 1.092983794203763
$ws.Range("I5").Value = 1.050124853207787
$ws.Range("J5").Value = 1.092171103645076
$ws.Range("K5").Value = 1.076101467698063
$ws.Range("L5").Value = 1.090089239482168
$ws.Range("M5").Value = 1.09528289391573
$ws.Range("N5").Value = 1.093722111886029
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.088086292395382
$ws.Range("D6").Value = 1.073819302758494
$ws.Range("E6").Value = 1.087879524901573
$ws.Range("F6").Value = 1.093074079573157
$ws.Range("I6").Value = 1.050145097602162
$ws.Range("J6").Value = 1.092274865630796
$ws.Range("K6").Value = 1.07615310433086
$ws.Range("L6").Value = 1.090181687529452
$ws.Range("M6").Value = 1.095364773330405
$ws.Range("N6").Value = 1.09382602122567
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.087263692960082
$ws.Range("D7").Value = 1.073405314510622
$ws.Range("E7").Value = 1.087184985357348
$ws.Range("F7").Value = 1.092452714507599
$ws.Range("I7").Value = 1.050005559312791
$ws.Range("J7").Value = 1.091560754380985
$ws.Range("K7").Value = 1.075797560846316
$ws.Range("L7").Value = 1.089545390966932
$ws.Range("M7").Value = 1.094801175209114
$ws.Range("N7").Value = 1.093110895855974
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.083809499867155
$ws.Range("D8").Value = 1.071663900687967
$ws.Range("E8").Value = 1.084267271472018
$ws.Range("F8").Value = 1.089841667507485
$ws.Range("I8").Value = 1.049413926702445
$ws.Range("J8").Value = 1.088560023864798
$ws.Range("K8").Value = 1.074299343277429
$ws.Range("L8").Value = 1.086870392364584
$ws.Range("M8").Value = 1.092430749599932
$ws.Range("N8").Value = 1.090105903958266
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.077659494217236
$ws.Range("D9").Value = 1.068552373537535
$ws.Range("E9").Value = 1.079067792657506
$ws.Range("F9").Value = 1.085186087343976
$ws.Range("I9").Value = 1.048339153687703
$ws.Range("J9").Value = 1.08320951431182
$ws.Range("K9").Value = 1.071612294086273
$ws.Range("L9").Value = 1.082095994720611
$ws.Range("M9").Value = 1.088196142281816
$ws.Range("N9").Value = 1.084747796068011
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.073513620801457
$ws.Range("D10").Value = 1.066447785812249
$ws.Range("E10").Value = 1.07555968074041
$ws.Range("F10").Value = 1.082043373217774
$ws.Range("I10").Value = 1.047600353985417
$ws.Range("J10").Value = 1.079597349786759
$ws.Range("K10").Value = 1.069788012553208
$ws.Range("L10").Value = 1.078869646635455
$ws.Range("M10").Value = 1.085332146742193
$ws.Range("N10").Value = 1.081130501855004
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.071706780382016
$ws.Range("D11").Value = 1.065528993265787
$ws.Range("E11").Value = 1.074030098949303
$ws.Range("F11").Value = 1.08067277464093
$ws.Range("I11").Value = 1.047275021668271
$ws.Range("J11").Value = 1.078021874149623
$ws.Range("K11").Value = 1.068989972439615
$ws.Range("L11").Value = 1.077461716768334
$ws.Range("M11").Value = 1.084081803171332
$ws.Range("N11").Value = 1.079552788861753
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.071033823676628
$ws.Range("D12").Value = 1.065186559949532
$ws.Range("E12").Value = 1.073460305892511
$ws.Range("F12").Value = 1.080162158230134
$ws.Range("I12").Value = 1.04715335123192
$ws.Range("J12").Value = 1.077434904659013
$ws.Range("K12").Value = 1.068692299565211
$ws.Range("L12").Value = 1.07693706012955
$ws.Range("M12").Value = 1.083615792047071
$ws.Range("N12").Value = 1.07896498580714
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.071178258407232
$ws.Range("D13").Value = 1.065260065672097
$ws.Range("E13").Value = 1.073582603450609
$ws.Range("F13").Value = 1.080271756483508
$ws.Range("I13").Value = 1.047179487599514
$ws.Range("J13").Value = 1.077560892533283
$ws.Range("K13").Value = 1.068756208184535
$ws.Range("L13").Value = 1.077049678016479
$ws.Range("M13").Value = 1.083715825121287
$ws.Range("N13").Value = 1.079091152598642
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.071651190994977
$ws.Range("D14").Value = 1.065500711298847
$ws.Range("E14").Value = 1.073983033413487
$ws.Range("F14").Value = 1.080630598102167
$ws.Range("I14").Value = 1.047264981304781
$ws.Range("J14").Value = 1.077973391464789
$ws.Range("K14").Value = 1.068965392263242
$ws.Range("L14").Value = 1.077418383238263
$ws.Range("M14").Value = 1.08404331498767
$ws.Range("N14").Value = 1.079504237325945
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.071942337787875
$ws.Range("D15").Value = 1.065648827421854
$ws.Range("E15").Value = 1.074229532812346
$ws.Range("F15").Value = 1.080851490131546
$ws.Range("I15").Value = 1.047317546774142
$ws.Range("J15").Value = 1.078227309699833
$ws.Range("K15").Value = 1.069094111677872
$ws.Range("L15").Value = 1.077645329361315
$ws.Range("M15").Value = 1.084244881913363
$ws.Range("N15").Value = 1.079758516154008
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.073633283844276
$ws.Range("D16").Value = 1.066508602885047
$ws.Range("E16").Value = 1.075660967281767
$ws.Range("F16").Value = 1.082134125374707
$ws.Range("I16").Value = 1.047621829844446
$ws.Range("J16").Value = 1.079701664194411
$ws.Range("K16").Value = 1.069840802669922
$ws.Range("L16").Value = 1.078962852309555
$ws.Range("M16").Value = 1.085414909049393
$ws.Range("N16").Value = 1.081234964401079
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.074690804853973
$ws.Range("D17").Value = 1.067045892413368
$ws.Range("E17").Value = 1.076556007150682
$ws.Range("F17").Value = 1.082936037786383
$ws.Range("I17").Value = 1.047811237074969
$ws.Range("J17").Value = 1.080623399786254
$ws.Range("K17").Value = 1.07030699151216
$ws.Range("L17").Value = 1.079786345630038
$ws.Range("M17").Value = 1.086146071203217
$ws.Range("N17").Value = 1.082158008963187
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.075306519411961
$ws.Range("D18").Value = 1.067358563289652
$ws.Range("E18").Value = 1.077077054900718
$ws.Range("F18").Value = 1.083402838958235
$ws.Range("I18").Value = 1.047921192300319
$ws.Range("J18").Value = 1.081159938045723
$ws.Range("K18").Value = 1.070578130474674
$ws.Range("L18").Value = 1.080265627185463
$ws.Range("M18").Value = 1.086571563121416
$ws.Range("N18").Value = 1.082695309168526
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.075516274232273
$ws.Range("D19").Value = 1.067465054497893
$ws.Range("E19").Value = 1.077254548364346
$ws.Range("F19").Value = 1.083561847883336
$ws.Range("I19").Value = 1.047958595877514
$ws.Range("J19").Value = 1.081342699796094
$ws.Range("K19").Value = 1.070670450150104
$ws.Range("L19").Value = 1.080428873559197
$ws.Range("M19").Value = 1.086716479500502
$ws.Range("N19").Value = 1.082878330461546
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.074577459127966
$ws.Range("D20").Value = 1.066988321105085
$ws.Range("E20").Value = 1.07646008309778
$ws.Range("F20").Value = 1.082850097807199
$ws.Range("I20").Value = 1.047790969635226
$ws.Range("J20").Value = 1.080524619902117
$ws.Range("K20").Value = 1.070257054842013
$ws.Range("L20").Value = 1.07969810137409
$ws.Range("M20").Value = 1.086067726368066
$ws.Range("N20").Value = 1.082059088800284
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.071511974790793
$ws.Range("D21").Value = 1.065429879205918
$ws.Range("E21").Value = 1.073865162432136
$ws.Range("F21").Value = 1.080524970381733
$ws.Range("I21").Value = 1.047239828482309
$ws.Range("J21").Value = 1.077851970073665
$ws.Range("K21").Value = 1.068903827353802
$ws.Range("L21").Value = 1.077309855775896
$ws.Range("M21").Value = 1.083946921292564
$ws.Range("N21").Value = 1.079382643502519
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.069574044467936
$ws.Range("D22").Value = 1.064443340926646
$ws.Range("E22").Value = 1.07222412677687
$ws.Range("F22").Value = 1.079054281570091
$ws.Range("I22").Value = 1.04688851011299
$ws.Range("J22").Value = 1.076161312476827
$ws.Range("K22").Value = 1.068045781530429
$ws.Range("L22").Value = 1.075798473958834
$ws.Range("M22").Value = 1.082604335923974
$ws.Range("N22").Value = 1.077689584977973
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.070602399163525
$ws.Range("D23").Value = 1.064966966357102
$ws.Range("E23").Value = 1.073094990161873
$ws.Range("F23").Value = 1.079834770019478
$ws.Range("I23").Value = 1.047075209182236
$ws.Range("J23").Value = 1.07705855370947
$ws.Range("K23").Value = 1.06850134093537
$ws.Range("L23").Value = 1.076600632086703
$ws.Range("M23").Value = 1.08331694792051
$ws.Range("N23").Value = 1.07858810039608
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.074628678638875
$ws.Range("D24").Value = 1.067014337327994
$ws.Range("E24").Value = 1.07650343017856
$ws.Range("F24").Value = 1.082888933289469
$ws.Range("I24").Value = 1.047800129233805
$ws.Range("J24").Value = 1.08056925765912
$ws.Range("K24").Value = 1.07027962148624
$ws.Range("L24").Value = 1.079737978391303
$ws.Range("M24").Value = 1.086103130060988
$ws.Range("N24").Value = 1.082103789948022
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.079257255010217
$ws.Range("D25").Value = 1.06936200494748
$ws.Range("E25").Value = 1.080419145084964
$ws.Range("F25").Value = 1.086396368141518
$ws.Range("I25").Value = 1.048620890646107
$ws.Range("J25").Value = 1.084600491905127
$ws.Range("K25").Value = 1.072312663994931
$ws.Range("L25").Value = 1.083337749620451
$ws.Range("M25").Value = 1.089297939791447
$ws.Range("N25").Value = 1.086140749009048
